# guardrail / 3.xlsx — remove the "skip empty row" guard so the computed
# region can legitimately contain blank rows (row 19 stays empty on purpose,
# splitting the K:O block into rows 11:18 and rows 20:27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing data tweak ---------------------------------------------------
$ws.Range("K8").Value = 2222222

# --- first new block: rows 11-18 (row 19 intentionally left empty) --------
$ws.Range("K11").Value = 3123
$ws.Range("L11").Value = 100
$ws.Range("M11").Value = 130
$ws.Range("N11").Value = 20

$ws.Range("K12").Value = 3123
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = 130
$ws.Range("N12").Value = 20

$ws.Range("K13").Value = 3123
$ws.Range("L13").Value = 100
$ws.Range("M13").Value = 130
$ws.Range("N13").Value = 20

$ws.Range("K14").Value = 322456
$ws.Range("L14").Value = 120
$ws.Range("M14").Value = 110
$ws.Range("N14").Value = 40

$ws.Range("K15").Value = 3456
$ws.Range("L15").Value = 120
$ws.Range("M15").Value = 110
$ws.Range("N15").Value = 40

$ws.Range("K16").Value = 111
$ws.Range("L16").Value = -60
$ws.Range("M16").Value = -10
$ws.Range("N16").Value = 30

$ws.Range("K17").Value = 2323
$ws.Range("L17").Value = -60
$ws.Range("M17").Value = -10
$ws.Range("N17").Value = 30

$ws.Range("K18").Value = 3969
$ws.Range("L18").Value = -70
$ws.Range("M18").Value = -60
$ws.Range("N18").Value = 30

# One shared formula covering O11:O18 (mirrors the O2:O9 block above it).
$ws.Range("O11:O18").Formula = "=M11-L11"

# Manual "in tolerance" highlight (yellow fill == existing style index 1).
$ws.Range("O11").Interior.Color = 65535
$ws.Range("O12").Interior.Color = 65535
$ws.Range("O13").Interior.Color = 65535
$ws.Range("O16").Interior.Color = 65535
$ws.Range("O17").Interior.Color = 65535

# --- row 19 is left completely blank on purpose ----------------------------

# --- second new block: rows 20-27 ------------------------------------------
$ws.Range("K20").Value = 3123
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = 130
$ws.Range("N20").Value = 20

$ws.Range("K21").Value = 3123
$ws.Range("L21").Value = 100
$ws.Range("M21").Value = 130
$ws.Range("N21").Value = 20

$ws.Range("K22").Value = 1111
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 130
$ws.Range("N22").Value = 20

$ws.Range("K23").Value = 3456
$ws.Range("L23").Value = 120
$ws.Range("M23").Value = 110
$ws.Range("N23").Value = 40

$ws.Range("K24").Value = 3456
$ws.Range("L24").Value = 120
$ws.Range("M24").Value = 110
$ws.Range("N24").Value = 40

$ws.Range("K25").Value = 444444444
$ws.Range("L25").Value = -60
$ws.Range("M25").Value = -10
$ws.Range("N25").Value = 30

$ws.Range("K26").Value = 33333333333
$ws.Range("L26").Value = -60
$ws.Range("M26").Value = -10
$ws.Range("N26").Value = 30

$ws.Range("K27").Value = 3969
$ws.Range("L27").Value = -70
$ws.Range("M27").Value = -60
$ws.Range("N27").Value = 30

$ws.Range("O20:O27").Formula = "=M20-L20"

$ws.Range("O20").Interior.Color = 65535
$ws.Range("O21").Interior.Color = 65535
$ws.Range("O22").Interior.Color = 65535
$ws.Range("O25").Interior.Color = 65535
$ws.Range("O26").Interior.Color = 65535

# --- column K width (bestFit-style autosize, now that long values exist) --
$ws.Columns("K").ColumnWidth = 11.1666666666667

# --- view state: scroll the visible pane over a bit, select K35 -----------
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$ws.Range("K35").Select()

# --- page setup (kept as portrait / paper size 9 = A4, matching target) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
